$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header values) ---
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# --- Row 2 ---
$ws.Range("B2").Value = 10.43629856256954
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = 13.983032143451659
$ws.Range("E2").ClearContents()

# --- Row 3 ---
$ws.Range("B3").Value = 10.228377126803458
$ws.Range("C3").ClearContents()
$ws.Range("D3").Value = 12.070514308248143
$ws.Range("E3").Value = -5.1271008514283096

# --- Update selection to match the new used range of interest ---
$ws.Range("B1:E3").Select()
